$wb = $excel.ActiveWorkbook

# --- "Modify Transaction" sheet: insert a new row for ReceiptNumber ---
$ws1 = $wb.Worksheets.Item("Modify Transaction")

# Insert a new row above the current row 4 (submitmakerepayment/click),
# pushing it down to row 5.
$ws1.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the ReceiptNumber label/value.
$ws1.Cells.Item(4, 1).Value = "ReceiptNumber"
$ws1.Cells.Item(4, 2).Value = 1234

# --- Selections / active sheet ---
# "Transactions" was previously the active sheet (activeTab=3) with
# selection D6. Make "Modify Transaction" the active sheet instead,
# with selection E8, leaving "Transactions" selection untouched at D6.
$ws4 = $wb.Worksheets.Item("Transactions")
$ws4.Range("D6").Select()

$ws1.Activate()
$ws1.Range("E8").Select()
